$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2793.9185
$ws.Range("I137").Value = 2600.0278
$ws.Range("J137").Value = 3330.8462
$ws.Range("K137").Value = 7800.0834
$ws.Range("L137").Value = 9992.5386
$ws.Range("M137").Value = -5250.0834
$ws.Range("N137").Value = -15092.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 22705.428
$ws.Range("J37").Value = 23986.834
$ws.Range("L37").Value = 23986.834
$ws.Range("N37").Value = -24532.834

$ws.Range("H74").Value = 5352.8076
$ws.Range("I74").Value = 5637.95
$ws.Range("J74").Value = 4402.3335
$ws.Range("K74").Value = 5637.95
$ws.Range("L74").Value = 4402.3335
$ws.Range("M74").Value = -4763.95
$ws.Range("N74").Value = -6150.3335

$ws.Range("H77").Value = 5352.8076
$ws.Range("I77").Value = 5637.95
$ws.Range("J77").Value = 4402.3335
$ws.Range("K77").Value = 28189.75
$ws.Range("L77").Value = 22011.6675
$ws.Range("M77").Value = -23821.75
$ws.Range("N77").Value = -30747.6675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 885.5
$ws.Range("I107").Value = 818.55
$ws.Range("J107").Value = 1019.4
$ws.Range("K107").Value = 818.55
$ws.Range("L107").Value = 1019.4
$ws.Range("M107").Value = 1101.45
$ws.Range("N107").Value = -4859.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1065.2142
$ws.Range("I22").Value = 1238.5
$ws.Range("J22").Value = 632
$ws.Range("K22").Value = 1238.5
$ws.Range("L22").Value = 632
$ws.Range("M22").Value = -888.5
$ws.Range("N22").Value = -1332

$ws.Range("H31").Value = 28871.176
$ws.Range("I31").Value = 39972.82
$ws.Range("J31").Value = 2967.3333
$ws.Range("K31").Value = 39972.82
$ws.Range("L31").Value = 2967.3333
$ws.Range("M31").Value = -39677.82
$ws.Range("N31").Value = -3557.3333

$ws.Range("H34").Value = 28871.176
$ws.Range("I34").Value = 39972.82
$ws.Range("J34").Value = 2967.3333
$ws.Range("K34").Value = 39972.82
$ws.Range("L34").Value = 2967.3333
$ws.Range("M34").Value = -39770.82
$ws.Range("N34").Value = -3371.3333

$ws.Range("H35").Value = 19049.5
$ws.Range("I35").Value = 3574.25
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 3574.25
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -3280.25
$ws.Range("N35").Value = -50588

$ws.Range("H58").Value = 1105.2
$ws.Range("I58").Value = 1003
$ws.Range("J58").Value = 1514
$ws.Range("K58").Value = 1003
$ws.Range("L58").Value = 1514
$ws.Range("M58").Value = -800
$ws.Range("N58").Value = -1920

$ws.Range("H132").Value = 14924.75
$ws.Range("I132").Value = 10618
$ws.Range("J132").Value = 24399.6
$ws.Range("K132").Value = 31854
$ws.Range("L132").Value = 73198.79999999999
$ws.Range("M132").Value = -29324
$ws.Range("N132").Value = -78258.79999999999

$ws.Range("H135").Value = 45045.57
$ws.Range("J135").Value = 45045.57
$ws.Range("L135").Value = 45045.57
$ws.Range("N135").Value = -55185.57

$ws.Range("H136").Value = 1105.2
$ws.Range("I136").Value = 1003
$ws.Range("J136").Value = 1514
$ws.Range("K136").Value = 3009
$ws.Range("L136").Value = 4542
$ws.Range("M136").Value = -459
$ws.Range("N136").Value = -9642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 666
$ws.Range("I5").Value = 400.85184
$ws.Range("J5").Value = 1316.8182
$ws.Range("K5").Value = 1202.55552
$ws.Range("L5").Value = 3950.4546
$ws.Range("M5").Value = -1090.55552
$ws.Range("N5").Value = -4174.4546

$ws.Range("H48").Value = 37038530
$ws.Range("I48").Value = 475
$ws.Range("J48").Value = 41668284
$ws.Range("K48").Value = 1425
$ws.Range("L48").Value = 125004852
$ws.Range("M48").Value = -1175
$ws.Range("N48").Value = -125005352

$ws.Range("H55").Value = 3208.3333
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 3363.6365
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 10090.9095
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -10444.9095

$ws.Range("H104").Value = 5998.5557
$ws.Range("J104").Value = 6498.375
$ws.Range("L104").Value = 19495.125
$ws.Range("N104").Value = -24737.125

$ws.Range("H105").Value = 9346.666999999999
$ws.Range("J105").Value = 9346.666999999999
$ws.Range("L105").Value = 28040.001
$ws.Range("N105").Value = -33282.001

$ws.Range("H106").Value = 3402.9
$ws.Range("J106").Value = 3402.9
$ws.Range("L106").Value = 10208.7
$ws.Range("N106").Value = -12100.7

$ws.Range("H118").Value = 1948.25
$ws.Range("I118").Value = 403.375
$ws.Range("J118").Value = 2566.2
$ws.Range("K118").Value = 1210.125
$ws.Range("L118").Value = 7698.599999999999
$ws.Range("M118").Value = 32.875
$ws.Range("N118").Value = -10184.6

$ws.Range("H120").Value = 21346.4
$ws.Range("I120").Value = 2850
$ws.Range("J120").Value = 33677.332
$ws.Range("K120").Value = 8550
$ws.Range("L120").Value = 101031.996
$ws.Range("M120").Value = -3712
$ws.Range("N120").Value = -110707.996

$ws.Range("H121").Value = 6166.6665
$ws.Range("I121").Value = 500
$ws.Range("K121").Value = 1500
$ws.Range("M121").Value = -190

$ws.Range("H134").Value = 3106.0386
$ws.Range("I134").Value = 1492.579
$ws.Range("J134").Value = 7485.4287
$ws.Range("K134").Value = 4477.737
$ws.Range("L134").Value = 22456.2861
$ws.Range("M134").Value = 592.2629999999999
$ws.Range("N134").Value = -32596.2861

$ws.Range("H135").Value = 666
$ws.Range("I135").Value = 400.85184
$ws.Range("J135").Value = 1316.8182
$ws.Range("K135").Value = 3607.66656
$ws.Range("L135").Value = 11851.3638
$ws.Range("M135").Value = -1072.66656
$ws.Range("N135").Value = -16921.3638

$ws.Range("H138").Value = 2408.2
$ws.Range("I138").Value = 887.25
$ws.Range("J138").Value = 3422.1667
$ws.Range("K138").Value = 2661.75
$ws.Range("L138").Value = 10266.5001
$ws.Range("M138").Value = 2478.25
$ws.Range("N138").Value = -20546.5001

$ws.Range("H141").Value = 2613.8462
$ws.Range("I141").Value = 1361.8182
$ws.Range("J141").Value = 9500
$ws.Range("K141").Value = 4085.4546
$ws.Range("L141").Value = 28500
$ws.Range("M141").Value = 1094.5454
$ws.Range("N141").Value = -38860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 24.066668
$ws.Range("I2").Value = 31
$ws.Range("J2").Value = 16.142857
$ws.Range("K2").Value = 31
$ws.Range("L2").Value = 16.142857
$ws.Range("M2").Value = 82
$ws.Range("N2").Value = -242.142857

$ws.Range("H122").Value = 2112.5
$ws.Range("I122").Value = 1620
$ws.Range("J122").Value = 2464.2856
$ws.Range("K122").Value = 4860
$ws.Range("L122").Value = 7392.8568
$ws.Range("M122").Value = -2410
$ws.Range("N122").Value = -12292.8568

$ws.Range("H134").Value = 11108.667
$ws.Range("J134").Value = 11108.667
$ws.Range("L134").Value = 33326.001
$ws.Range("N134").Value = -38396.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 591.46155
$ws.Range("I22").Value = 498.9
$ws.Range("K22").Value = 498.9
$ws.Range("M22").Value = -203.9

$ws.Range("H27").Value = 591.46155
$ws.Range("I27").Value = 498.9
$ws.Range("K27").Value = 498.9
$ws.Range("M27").Value = -391.9

$ws.Range("H32").Value = 19432.223
$ws.Range("I32").Value = 976.6667
$ws.Range("J32").Value = 56343.332
$ws.Range("K32").Value = 976.6667
$ws.Range("L32").Value = 56343.332
$ws.Range("M32").Value = -659.6667
$ws.Range("N32").Value = -56977.332

$ws.Range("H122").Value = 3939.9
$ws.Range("I122").Value = 4282.2354
$ws.Range("K122").Value = 12846.7062
$ws.Range("M122").Value = -10396.7062

$ws.Range("H132").Value = 58927.39
$ws.Range("I132").Value = 73942.14
$ws.Range("J132").Value = 6375.75
$ws.Range("K132").Value = 221826.42
$ws.Range("L132").Value = 19127.25
$ws.Range("M132").Value = -219296.42
$ws.Range("N132").Value = -24187.25
